# Apply the workbook edits described by the diff:
#  1. Update the report date in the merged title cell (A2) on every sheet
#     from "05-10-2018" to "06-01-2026". The cell is backed by a shared
#     string used by all three worksheets, so this single textual change
#     propagates everywhere the string is used.
#  2. Recolor the bold "Report on managers:" header font (A4 on every
#     sheet) from RGB(59,59,56) / #3B3B38 to RGB(75,70,42) / #4B462A.

$wb = $excel.ActiveWorkbook

$newHeaderColor = [Convert]::ToInt32("2A464B", 16)  # BGR-packed #4B462A (R=4B,G=46,B=2A)

foreach ($ws in $wb.Worksheets) {
    $ws.Range("A2").Value = "Date: 06-01-2026 - Department: Sales department"
    $ws.Range("A4").Font.Color = $newHeaderColor
}
